$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(13, 1).Value = "20160405_130217"
$ws.Cells.Item(13, 2).Value = 547.1
$ws.Cells.Item(13, 3).Value = "remove multiple spaces, trim `"space`" and `",`", convert unicode to ascii, convert to lower"
$ws.Cells.Item(13, 4).Value = "8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit"
$ws.Cells.Item(13, 5).Value = "Neuron Network"
$ws.Cells.Item(13, 6).Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Cells.Item(13, 7).Value = 0.993333333333333
$ws.Cells.Item(13, 8).Value = 0.99009900990099
$ws.Cells.Item(13, 9).Value = "0 filters: "
$ws.Cells.Item(13, 10).Value = 0.427083333333333

$ws.Cells.Item(14, 1).Value = "20160405_131125"
$ws.Cells.Item(14, 2).Value = 571.652
$ws.Cells.Item(14, 3).Value = "remove multiple spaces, trim `"space`" and `",`", convert unicode to ascii, convert to lower"
$ws.Cells.Item(14, 4).Value = "8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit"
$ws.Cells.Item(14, 5).Value = "Neuron Network"
$ws.Cells.Item(14, 6).Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Cells.Item(14, 7).Value = 0.99
$ws.Cells.Item(14, 8).Value = 0.99009900990099
$ws.Cells.Item(14, 9).Value = "0 filters: "
$ws.Cells.Item(14, 10).Value = 0.427083333333333

$ws.Cells.Item(15, 1).Value = "20160405_132056"
$ws.Cells.Item(15, 2).Value = 583.173
$ws.Cells.Item(15, 3).Value = "remove multiple spaces, trim `"space`" and `",`", convert unicode to ascii, convert to lower"
$ws.Cells.Item(15, 4).Value = "8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit"
$ws.Cells.Item(15, 5).Value = "Neuron Network"
$ws.Cells.Item(15, 6).Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Cells.Item(15, 7).Value = 0.990666666666667
$ws.Cells.Item(15, 8).Value = 0.99009900990099
$ws.Cells.Item(15, 9).Value = "0 filters: "
$ws.Cells.Item(15, 10).Value = 0.40625

$ws.Cells.Item(16, 1).Value = "20160405_133039"
$ws.Cells.Item(16, 2).Value = 587.907
$ws.Cells.Item(16, 3).Value = "remove multiple spaces, trim `"space`" and `",`", convert unicode to ascii, convert to lower"
$ws.Cells.Item(16, 4).Value = "8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit"
$ws.Cells.Item(16, 5).Value = "Neuron Network"
$ws.Cells.Item(16, 6).Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Cells.Item(16, 7).Value = 0.994
$ws.Cells.Item(16, 8).Value = 0.99009900990099
$ws.Cells.Item(16, 9).Value = "0 filters: "
$ws.Cells.Item(16, 10).Value = 0.416666666666667

$ws.Cells.Item(17, 1).Value = "20160405_134027"
$ws.Cells.Item(17, 2).Value = 597.983
$ws.Cells.Item(17, 3).Value = "remove multiple spaces, trim `"space`" and `",`", convert unicode to ascii, convert to lower"
$ws.Cells.Item(17, 4).Value = "8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit"
$ws.Cells.Item(17, 5).Value = "Neuron Network"
$ws.Cells.Item(17, 6).Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Cells.Item(17, 7).Value = 0.991333333333333
$ws.Cells.Item(17, 8).Value = 0.99009900990099
$ws.Cells.Item(17, 9).Value = "0 filters: "
$ws.Cells.Item(17, 10).Value = 0.416666666666667

$ws.Cells.Item(18, 1).Value = "20160405_145007"
$ws.Cells.Item(18, 2).Value = 1093.521
$ws.Cells.Item(18, 3).Value = "convert unicode to ascii, convert to lower, remove multiple spaces, trim `"space`" and `",`""
$ws.Cells.Item(18, 4).Value = "8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit"
$ws.Cells.Item(18, 5).Value = "Neuron Network"
$ws.Cells.Item(18, 6).Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Cells.Item(18, 7).Value = 0.990666666666667
$ws.Cells.Item(18, 8).Value = 0.986798679867987
$ws.Cells.Item(18, 9).Value = "0 filters: "
$ws.Cells.Item(18, 10).Value = 0.326315789473684

$ws.Cells.Item(19, 1).Value = "20160405_150820"
$ws.Cells.Item(19, 2).Value = 1116.069
$ws.Cells.Item(19, 3).Value = "convert unicode to ascii, convert to lower, remove multiple spaces, trim `"space`" and `",`""
$ws.Cells.Item(19, 4).Value = "8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit"
$ws.Cells.Item(19, 5).Value = "Neuron Network"
$ws.Cells.Item(19, 6).Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Cells.Item(19, 7).Value = 0.99
$ws.Cells.Item(19, 8).Value = 0.99009900990099
$ws.Cells.Item(19, 9).Value = "0 filters: "
$ws.Cells.Item(19, 10).Value = 0.416666666666667

$ws.Cells.Item(20, 1).Value = "20160405_152656"
$ws.Cells.Item(20, 2).Value = 1162.183
$ws.Cells.Item(20, 3).Value = "convert unicode to ascii, convert to lower, remove multiple spaces, trim `"space`" and `",`""
$ws.Cells.Item(20, 4).Value = "8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit"
$ws.Cells.Item(20, 5).Value = "Neuron Network"
$ws.Cells.Item(20, 6).Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Cells.Item(20, 7).Value = 0.992666666666667
$ws.Cells.Item(20, 8).Value = 0.99009900990099
$ws.Cells.Item(20, 9).Value = "0 filters: "
$ws.Cells.Item(20, 10).Value = 0.416666666666667

$ws.Cells.Item(21, 1).Value = "20160405_154619"
$ws.Cells.Item(21, 2).Value = 1186.798
$ws.Cells.Item(21, 3).Value = "convert unicode to ascii, convert to lower, remove multiple spaces, trim `"space`" and `",`""
$ws.Cells.Item(21, 4).Value = "8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit"
$ws.Cells.Item(21, 5).Value = "Neuron Network"
$ws.Cells.Item(21, 6).Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Cells.Item(21, 7).Value = 0.996666666666667
$ws.Cells.Item(21, 8).Value = 0.99009900990099
$ws.Cells.Item(21, 9).Value = "0 filters: "
$ws.Cells.Item(21, 10).Value = 0.416666666666667

$ws.Cells.Item(22, 1).Value = "20160405_160605"
$ws.Cells.Item(22, 2).Value = 1218.147
$ws.Cells.Item(22, 3).Value = "convert unicode to ascii, convert to lower, remove multiple spaces, trim `"space`" and `",`""
$ws.Cells.Item(22, 4).Value = "8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit"
$ws.Cells.Item(22, 5).Value = "Neuron Network"
$ws.Cells.Item(22, 6).Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Cells.Item(22, 7).Value = 0.996
$ws.Cells.Item(22, 8).Value = 0.99009900990099
$ws.Cells.Item(22, 9).Value = "0 filters: "
$ws.Cells.Item(22, 10).Value = 0.416666666666667

